$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season-record columns, styled like the
# existing header row (bold, centered, thin border - same as AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
